# Update cryptocurrency price/volume figures per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainValue($range, $value) {
    $range.Value = $value
}

function Set-TextValue($range, $value) {
    # Force the cell to keep its value as literal text (preserving trailing
    # zeros / exact digit count) instead of letting Excel auto-convert a
    # numeric-looking string into a number. Style is restored afterwards so
    # the cell's formatting is left untouched.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-PlainValue $ws.Range("D2") "30.483.98"
Set-PlainValue $ws.Range("E2") "  +0.24%  "
Set-PlainValue $ws.Range("D3") "1.913.99"
Set-PlainValue $ws.Range("E3") "  -0.22%  "
Set-PlainValue $ws.Range("E4") "  +0.13%  "
Set-TextValue $ws.Range("D5") "245.00"
Set-PlainValue $ws.Range("E5") "  +1.54%  "
Set-TextValue $ws.Range("D6") "0.9998"
Set-PlainValue $ws.Range("E6") "  +0.07%  "
Set-TextValue $ws.Range("D7") "0.4835"
Set-PlainValue $ws.Range("E7") "  +2.84%  "
Set-TextValue $ws.Range("D8") "0.2892"
Set-PlainValue $ws.Range("E8") "  +1.35%  "
Set-TextValue $ws.Range("D9") "0.06717"
Set-PlainValue $ws.Range("E9") "  -1.71%  "
Set-TextValue $ws.Range("D10") "109.96"
Set-PlainValue $ws.Range("E10") "  +0.40%  "
Set-TextValue $ws.Range("D11") "19.01"
Set-PlainValue $ws.Range("E11") "  +3.90%  "
Set-PlainValue $ws.Range("D12") "1.916.44"
Set-PlainValue $ws.Range("E12") "  -0.04%  "
Set-TextValue $ws.Range("D13") "0.07549"
Set-PlainValue $ws.Range("E13") "  -1.12%  "
Set-TextValue $ws.Range("D14") "5.271"
Set-PlainValue $ws.Range("E14") "  +1.32%  "
Set-TextValue $ws.Range("D15") "0.6717"
Set-PlainValue $ws.Range("E15") "  +2.14%  "
Set-TextValue $ws.Range("D16") "282.22"
Set-PlainValue $ws.Range("E16") "  -3.26%  "
Set-PlainValue $ws.Range("D17") "30.498.33"
Set-PlainValue $ws.Range("E17") "  +0.24%  "
Set-TextValue $ws.Range("D18") "0.9998"
Set-PlainValue $ws.Range("E18") "  -0.01%  "
Set-TextValue $ws.Range("D19") "0.000007567"
Set-PlainValue $ws.Range("E19") "  -1.09%  "
Set-PlainValue $ws.Range("E20") "  -0.47%  "
Set-PlainValue $ws.Range("D21") "2.167.13"
Set-PlainValue $ws.Range("E21") "  +0.69%  "
Set-TextValue $ws.Range("D22") "5.489"
Set-PlainValue $ws.Range("E22") "  +4.73%  "
Set-TextValue $ws.Range("D23") "0.9997"
Set-PlainValue $ws.Range("E23") "  +0.12%  "
Set-TextValue $ws.Range("D24") "6.449"
Set-PlainValue $ws.Range("E24") "  +3.78%  "
Set-TextValue $ws.Range("D25") "9.469"
Set-PlainValue $ws.Range("E25") "  +2.06%  "
Set-TextValue $ws.Range("D26") "164.23"
Set-PlainValue $ws.Range("E26") "  -2.15%  "
Set-PlainValue $ws.Range("E27") "  -6.35%  "
Set-TextValue $ws.Range("D28") "2.123"
Set-PlainValue $ws.Range("E28") "  +3.69%  "
Set-TextValue $ws.Range("D29") "0.1053"
Set-PlainValue $ws.Range("E29") "  -1.81%  "
Set-TextValue $ws.Range("D30") "1.406"
Set-PlainValue $ws.Range("E30") "  +2.63%  "
Set-TextValue $ws.Range("D31") "4.151"
Set-PlainValue $ws.Range("E31") "  -0.11%  "
Set-TextValue $ws.Range("D32") "4.038"
Set-PlainValue $ws.Range("E32") "  +1.99%  "
Set-TextValue $ws.Range("D33") "0.04988"
Set-PlainValue $ws.Range("E33") "  -1.15%  "
Set-TextValue $ws.Range("D34") "0.7298"
Set-PlainValue $ws.Range("E34") "  -1.47%  "
Set-PlainValue $ws.Range("E35") "  -1.09%  "
Set-TextValue $ws.Range("D36") "0.9995"
Set-PlainValue $ws.Range("E36") "  +0.11%  "
Set-TextValue $ws.Range("D37") "2.724"
Set-PlainValue $ws.Range("E37") "  -0.77%  "
Set-PlainValue $ws.Range("E38") "  -1.25%  "
Set-PlainValue $ws.Range("E39") "  -0.72%  "
Set-TextValue $ws.Range("D40") "110.75"
Set-PlainValue $ws.Range("E40") "  +1.97%  "
Set-TextValue $ws.Range("D41") "2.017"
Set-PlainValue $ws.Range("E41") "  -1.82%  "
Set-TextValue $ws.Range("D42") "0.4449"
Set-PlainValue $ws.Range("E42") "  +5.20%  "
Set-TextValue $ws.Range("D43") "0.8648"
Set-PlainValue $ws.Range("E43") "  -1.39%  "
Set-TextValue $ws.Range("D44") "5.790"
Set-PlainValue $ws.Range("E44") "  -1.10%  "
Set-TextValue $ws.Range("D45") "0.9997"
Set-PlainValue $ws.Range("E45") "  +0.10%  "
Set-TextValue $ws.Range("D46") "67.95"
Set-PlainValue $ws.Range("E46") "  +0.29%  "
Set-TextValue $ws.Range("D47") "7.347"
Set-PlainValue $ws.Range("E47") "  +2.06%  "
Set-TextValue $ws.Range("D48") "49.17"
Set-PlainValue $ws.Range("E48") "  -5.42%  "
Set-TextValue $ws.Range("D49") "9.289"
Set-PlainValue $ws.Range("E49") "  +0.99%  "
Set-PlainValue $ws.Range("E50") "  +2.54%  "
Set-TextValue $ws.Range("D51") "34.82"
Set-PlainValue $ws.Range("E51") "  +0.03%  "